# Tenant page Add, Update & Delete sheets
# Adds three new worksheets (AddTenant, UpdateTenant, DeleteTenant) at the
# end of the workbook, populates the AddTenant sheet with sample data, and
# makes AddTenant the active sheet (matching the original author's edit).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---- AddTenant ---------------------------------------------------------
$sheetAdd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sheetAdd.Name = "AddTenant"

$sheetAdd.Range("A1").Value = "Tenant_name"
$sheetAdd.Range("B1").Value = "Product_code_search_input"
$sheetAdd.Range("C1").Value = "Assigned_count"
$sheetAdd.Range("D1").Value = "User_code_search_input"

$sheetAdd.Range("A2").Value = "tenant1"
$sheetAdd.Range("B2").Value = "TDC"
$sheetAdd.Range("C2").Value = 1
$sheetAdd.Range("D2").Value = "Test"

$sheetAdd.Range("A3").Value = "tenant2"
$sheetAdd.Range("B3").Value = "TDC"
$sheetAdd.Range("C3").Value = 1
$sheetAdd.Range("D3").Value = "Test"

# page setup (0.75"/0.75"/1"/1"/0.5"/0.5" -> points)
$sheetAdd.PageSetup.LeftMargin = 54
$sheetAdd.PageSetup.RightMargin = 54
$sheetAdd.PageSetup.TopMargin = 72
$sheetAdd.PageSetup.BottomMargin = 72
$sheetAdd.PageSetup.HeaderMargin = 36
$sheetAdd.PageSetup.FooterMargin = 36
$sheetAdd.PageSetup.LeftHeader = ""

$sheetAdd.Range("E11").Select() | Out-Null

# ---- UpdateTenant -------------------------------------------------------
$sheetUpdate = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetAdd)
$sheetUpdate.Name = "UpdateTenant"

$sheetUpdate.PageSetup.LeftMargin = 54
$sheetUpdate.PageSetup.RightMargin = 54
$sheetUpdate.PageSetup.TopMargin = 72
$sheetUpdate.PageSetup.BottomMargin = 72
$sheetUpdate.PageSetup.HeaderMargin = 36
$sheetUpdate.PageSetup.FooterMargin = 36
$sheetUpdate.PageSetup.LeftHeader = ""

$sheetUpdate.Range("A1").Select() | Out-Null

# ---- DeleteTenant ---------------------------------------------------------
$sheetDelete = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetUpdate)
$sheetDelete.Name = "DeleteTenant"

$sheetDelete.PageSetup.LeftMargin = 54
$sheetDelete.PageSetup.RightMargin = 54
$sheetDelete.PageSetup.TopMargin = 72
$sheetDelete.PageSetup.BottomMargin = 72
$sheetDelete.PageSetup.HeaderMargin = 36
$sheetDelete.PageSetup.FooterMargin = 36
$sheetDelete.PageSetup.LeftHeader = ""

$sheetDelete.Range("L23").Select() | Out-Null

# Make AddTenant the active (selected) sheet, as in the target workbook.
$sheetAdd.Activate()
